$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "test_subject"
$ws.Range("A3").Select() | Out-Null
